$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1777.8823
$ws.Range("J32").Value = 1887.5
$ws.Range("L32").Value = 1887.5
$ws.Range("N32").Value = -2539.5
$ws.Range("H62").Value = 7359304
$ws.Range("I62").Value = 12505217
$ws.Range("K62").Value = 12505217
$ws.Range("M62").Value = -12504593
$ws.Range("H65").Value = 7359304
$ws.Range("I65").Value = 12505217
$ws.Range("K65").Value = 62526085
$ws.Range("M65").Value = -62522965
$ws.Range("H137").Value = 3291.4644
$ws.Range("I137").Value = 2853.9443
$ws.Range("K137").Value = 8561.832900000001
$ws.Range("M137").Value = -6011.832900000001
$ws.Range("H138").Value = 6348.1064
$ws.Range("J138").Value = 7789.596
$ws.Range("L138").Value = 23368.788
$ws.Range("N138").Value = -33648.788

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2247.76
$ws.Range("I32").Value = 1540.6022
$ws.Range("K32").Value = 1540.6022
$ws.Range("M32").Value = -1253.6022
$ws.Range("H132").Value = 3306.8044
$ws.Range("I132").Value = 2740.8096
$ws.Range("K132").Value = 8222.4288
$ws.Range("M132").Value = -5692.4288

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 103999
$ws.Range("J59").Value = 103999
$ws.Range("L59").Value = 103999
$ws.Range("N59").Value = -105693
$ws.Range("H105").Value = 2323.76
$ws.Range("I105").Value = 2293.524
$ws.Range("K105").Value = 2293.524
$ws.Range("M105").Value = -546.5239999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2247.9092
$ws.Range("I58").Value = 1890.3334
$ws.Range("J58").Value = 3014.1428
$ws.Range("K58").Value = 1890.3334
$ws.Range("L58").Value = 3014.1428
$ws.Range("M58").Value = -1687.3334
$ws.Range("N58").Value = -3420.1428
$ws.Range("H99").Value = 5041.357
$ws.Range("J99").Value = 4681.2856
$ws.Range("L99").Value = 4681.2856
$ws.Range("N99").Value = -7677.2856
$ws.Range("H126").Value = 5041.357
$ws.Range("J126").Value = 4681.2856
$ws.Range("L126").Value = 14043.8568
$ws.Range("N126").Value = -18983.8568
$ws.Range("H132").Value = 2046.5758
$ws.Range("I132").Value = 1501.2759
$ws.Range("K132").Value = 4503.8277
$ws.Range("M132").Value = -1973.8277
$ws.Range("H136").Value = 2247.9092
$ws.Range("I136").Value = 1890.3334
$ws.Range("J136").Value = 3014.1428
$ws.Range("K136").Value = 5671.0002
$ws.Range("L136").Value = 9042.428400000001
$ws.Range("M136").Value = -3121.0002
$ws.Range("N136").Value = -14142.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 51767.668
$ws.Range("I109").Value = 2434.2
$ws.Range("K109").Value = 7302.599999999999
$ws.Range("M109").Value = -6262.599999999999
$ws.Range("H118").Value = 3308.0789
$ws.Range("J118").Value = 3451.484
$ws.Range("L118").Value = 10354.452
$ws.Range("N118").Value = -12840.452
$ws.Range("H131").Value = 3593.0908
$ws.Range("J131").Value = 3713.45
$ws.Range("L131").Value = 11140.35
$ws.Range("N131").Value = -21220.35

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 2014822
$ws.Range("J10").Value = 18527.5
$ws.Range("L10").Value = 18527.5
$ws.Range("N10").Value = -18865.5
$ws.Range("H18").Value = 500002500
$ws.Range("I18").Value = 5000
$ws.Range("J18").Value = 1000000000
$ws.Range("K18").Value = 5000
$ws.Range("L18").Value = 1000000000
$ws.Range("M18").Value = -4707
$ws.Range("N18").Value = -1000000586
$ws.Range("H126").Value = 3310.0557
$ws.Range("I126").Value = 2287.7778
$ws.Range("K126").Value = 6863.3334
$ws.Range("M126").Value = -4393.3334
$ws.Range("H132").Value = 64792.176
$ws.Range("I132").Value = 6341.6875
$ws.Range("K132").Value = 19025.0625
$ws.Range("M132").Value = -16495.0625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5594.2856
$ws.Range("I7").Value = 4965.909
$ws.Range("K7").Value = 4965.909
$ws.Range("M7").Value = -4853.909
$ws.Range("H48").Value = 8500
$ws.Range("I48").Value = 5000
$ws.Range("J48").Value = 12000
$ws.Range("K48").Value = 5000
$ws.Range("L48").Value = 12000
$ws.Range("M48").Value = -4339
$ws.Range("N48").Value = -13322
$ws.Range("H55").Value = 1741.091
$ws.Range("I55").Value = 374.5
$ws.Range("K55").Value = 374.5
$ws.Range("M55").Value = -201.5
$ws.Range("H95").Value = 40344
$ws.Range("J95").Value = 40344
$ws.Range("L95").Value = 40344
$ws.Range("N95").Value = -45836
$ws.Range("H122").Value = 720368.4
$ws.Range("I122").Value = 1670192.9
$ws.Range("K122").Value = 5010578.699999999
$ws.Range("M122").Value = -5008128.699999999
$ws.Range("H126").Value = 5594.2856
$ws.Range("I126").Value = 4965.909
$ws.Range("K126").Value = 14897.727
$ws.Range("M126").Value = -12427.727
$ws.Range("H136").Value = 367972.3
$ws.Range("I136").Value = 722694.4399999999
$ws.Range("K136").Value = 2168083.32
$ws.Range("M136").Value = -2165533.32

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8364.272000000001
$ws.Range("J62").Value = 7714.7144
$ws.Range("L62").Value = 7714.7144
$ws.Range("N62").Value = -8962.714400000001
$ws.Range("H65").Value = 8364.272000000001
$ws.Range("J65").Value = 7714.7144
$ws.Range("L65").Value = 38573.572
$ws.Range("N65").Value = -44813.572
$ws.Range("H69").Value = 30000
$ws.Range("J69").Value = 30000
$ws.Range("L69").Value = 30000
$ws.Range("N69").Value = -31498
$ws.Range("H72").Value = 30000
$ws.Range("J72").Value = 30000
$ws.Range("L72").Value = 90000
$ws.Range("N72").Value = -97488
$ws.Range("H126").Value = 1532.1154
$ws.Range("I126").Value = 1534.7916
$ws.Range("K126").Value = 4604.3748
$ws.Range("M126").Value = -2134.3748
$ws.Range("H132").Value = 46581.74
$ws.Range("I132").Value = 2169.5715
$ws.Range("J132").Value = 115667.336
$ws.Range("K132").Value = 6508.7145
$ws.Range("L132").Value = 347002.008
$ws.Range("M132").Value = -3978.7145
$ws.Range("N132").Value = -352062.008
$ws.Range("H133").Value = 199999
$ws.Range("J133").Value = 199999
$ws.Range("L133").Value = 199999
$ws.Range("N133").Value = -210119
$ws.Range("H136").Value = 806986.9
$ws.Range("I136").Value = 839925.3
$ws.Range("J136").Value = 675233
$ws.Range("K136").Value = 2519775.9
$ws.Range("L136").Value = 2025699
$ws.Range("M136").Value = -2517225.9
$ws.Range("N136").Value = -2030799
